# Quarterly indexing bug-fix: the date in column A was being written as the
# 1st of the first month of each quarter (Jan/Apr/Jul/Oct 1) instead of the
# intended mid-point "data available" date, which is the 15th of the
# *second* month of the quarter (one calendar month + 14 days later,
# e.g. 1988-07-01 -> 1988-08-15). Re-derive column A from itself using that
# corrected rule, for every data row (row 2 through the last used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date system epoch (serial 0 == 1899-12-30) so we can turn the
# numeric serial stored in the cell into a real date to do month math on.
$excelEpoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 150 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null -or $serial -eq "") { continue }

    $oldDate = $excelEpoch.AddDays([double]$serial)
    $shifted = $oldDate.AddMonths(1)
    $newDate = Get-Date -Year $shifted.Year -Month $shifted.Month -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate
}
